{"js": "const body = context.document.body;\n\n// 1. Remove the stray \"m.\" subdomain from the Facebook story link.\nconst linkResults = body.search(\"https://m.facebook.com/story.php\", { matchCase: true });\nlinkResults.load(\"items\");\nawait context.sync();\nif (linkResults.items.length > 0) {\n  linkResults.items[0].insertText(\"https://facebook.com/story.php\", Word.InsertLocation.replace);\n}\n\n// 2. Relocate the \"_GoBack\" bookmark (left behind by the last edit session)\n//    from the end of the Facebook-link paragraph to right after\n//    \"...Stay tuned for updates\" in the opening paragraph. Word only ever\n//    keeps a single \"_GoBack\" bookmark, so the old one is deleted first,\n//    exactly like the diff shows (bookmark removed from one spot, added\n//    at the other).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchorResults = body.search(\"Stay tuned for updates\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nconst anchor = anchorResults.items[0];\nconst afterAnchor = anchor.getRange(Word.RangeLocation.after);\nafterAnchor.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"m.\" subdomain from the Facebook story link.\n$find = $d.Content.Find\n$find.Text = \"https://m.facebook.com/story.php\"\n$find.Replacement.Text = \"https://facebook.com/story.php\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Relocate the \"_GoBack\" bookmark (left behind by the last edit session)\n#    from the end of the Facebook-link paragraph to right after\n#    \"...Stay tuned for updates\" in the opening paragraph. Word only ever\n#    keeps a single \"_GoBack\" bookmark, so re-adding it here removes the\n#    old one automatically, exactly like the diff shows.\n$rng = $d.Content\n$rng.Find.Text = \"Stay tuned for updates\"\n$rng.Find.Execute()\n$rng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
